$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.035.81'
$ws.Range("E2").Value = '  +1.08%  '
$ws.Range("D3").Value = '3.797.82'
$ws.Range("E3").Value = '  +2.35%  '
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '428.30'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +6.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.97'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +9.30%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.622'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +3.75%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.731'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +3.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.151'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -6.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000309'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -10.72%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.56'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +6.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.45'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +9.58%  '
$ws.Range("D14").Value = '4.382.92'
$ws.Range("E14").Value = '  +2.34%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.08'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +4.86%  '
$ws.Range("E16").Value = '  +0.54%  '
$ws.Range("D17").Value = '3.762.66'
$ws.Range("E17").Value = '  +1.69%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '19.90'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +3.09%  '
$ws.Range("E19").Value = '  +6.53%  '
$ws.Range("D20").Value = '66.119.76'
$ws.Range("E20").Value = '  +1.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '405.06'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.22%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.16'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +6.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.21'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +6.72%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.78'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.34%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '36.81'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +2.54%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.27'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +6.97%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.70'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +38.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.90'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +9.19%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.42'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.05%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '13.77'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +11.59%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '705.25'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +3.24%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.131'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +14.40%  '
$ws.Range("E33").Value = '  +3.23%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '41.02'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +7.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.19%  '
$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.65'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +34.48%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.148'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -3.41%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '56.57'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +2.70%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0473'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +5.40%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.64'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +41.32%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.86'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +3.73%  '
$ws.Range("D42").Value = '0.0₃0677'
$ws.Range("E42").Value = '  -5.47%  '
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.24%  '
$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.140'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +5.20%  '
$ws.Range("B45").Value = 'LidoDAOToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.38'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +5.72%  '
$ws.Range("B46").Value = 'TheGraph'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.322'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +12.22%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.16'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.42%  '
$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.65'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +4.44%  '
$ws.Range("B49").Value = 'ARBITRUM'
$ws.Range("C49").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.07'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +2.74%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '141.84'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.87%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.80'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.85%  '
